$d = $word.ActiveDocument
$anchor = $d.Paragraphs.Last
$baseIndex = $anchor.Index

# --- Phase 1: create all the new (empty) paragraphs first. Doing this
# up front -- before any text/formatting is applied -- means no paragraph
# can inherit stray "current formatting" (e.g. italic) left over from a
# sibling paragraph that is filled in later.
for ($k = 0; $k -lt 13; $k++) {
  $anchor.Range.InsertParagraphAfter()
  $anchor = $d.Paragraphs.Last
}

# --- Phase 2: fill in style/text/formatting for each paragraph. ---
# (Every paragraph created in Phase 1 is a clone of the anchor paragraph's
# style ("Title"), since each is freshly split off the still-Title-styled
# anchor -- so every one that is not Title needs an explicit style set.)
# --- paragraph 0: style=Heading1 ---
$para = $d.Paragraphs.Item($baseIndex + 1)
$para.Style = 'Heading1'
$pStart = $para.Range.Start
$para.Range.Text = 'Knärot – ekologi samt krav på livsmiljön'

# --- paragraph 1: style=Normal ---
$para = $d.Paragraphs.Item($baseIndex + 2)
$para.Style = 'Normal'
$pStart = $para.Range.Start
$para.Range.Text = 'Knärot är fridlyst enligt 8 och 15 §§ artskyddsförordningen och klassad som sårbar (VU) enligt rödlistan 2020. Knärot är beroende av hög och jämn luftfuktighet i gamla, ostörda skogsmiljöer och är känslig för snabba förändringar av ljus-/vindförhållanden eller uttorkning. På grund av ett alltför intensivt skogsbruk har den minskat med 40 (25-50) % under de senaste 60 åren och i framtiden bedöms minskningstakten uppgå till 30 (20-40) %. Till följd av att arten har en dokumenterat högre minskningstakt iförhållande till sin generationstid än vad som tidigare varit känt (data från Riksskogstaxeringen) höjdes den till hotkategori sårbar (VU) i rödlistan 2020 (Artdatabanken, 2021).'

# --- paragraph 2: style=Normal ---
$para = $d.Paragraphs.Item($baseIndex + 3)
$para.Style = 'Normal'
$pStart = $para.Range.Start
$para.Range.Text = 'Samuel Johnsons doktorsavhandling “Retention Forestry as a Conservation Measure for Boreal Forest Ground Vegetation“ (SLU, Uppsala 2014) visar att det krävs väl tilltagna skyddszoner för att knärotens växtplatser inte ska ta skada av skogsbruksåtgärder i intilliggande områden: “Study III shows that retention patches smaller than 0.5 ha do not lifeboat the sensitive forest herb G. repens, a species that depend on stable microclimatic conditions typical for intact forest stands.” Vidare “More sensitive forest species are not lifeboated in retention patches ranging from 0.05 to 0.5 ha (Papers II & III).”'
$d.Range($pStart + 34, $pStart + 116).Font.Italic = $true
$d.Range($pStart + 278, $pStart + 483).Font.Italic = $true
$d.Range($pStart + 490, $pStart + 608).Font.Italic = $true

# --- paragraph 3: style=Normal ---
$para = $d.Paragraphs.Item($baseIndex + 4)
$para.Style = 'Normal'
$pStart = $para.Range.Start
$para.Range.Text = 'Johnsons (2014) rekommendation på minst 50 meters breda skyddszoner runt knärotens växtplatser motsvarar en areal på 0,78 hektar, vilket ligger i linje med andra studier som gjorts på känsliga skogsarter: “In study III I also show that translocated specimens of G. repens survives well in mature forests at least 50 m from the nearest edge to an open area. Moreover, measures of temperature and humidity show that such distances from an open area is far enough to offer a microclimate that is more stable compared to what present in retention patches of around 0.1 ha. This means that the very centre of a circular patch with radius 50 m (equals a size of 0.78 ha) should offer conditions similar to interior forest and would perhaps be a suitable habitat for G. repens and similar species. Previous studies from both North America and Sweden have also concluded that patches between 0.5 and one ha are sufficient for preserving interior forest vegetation as well as sensitive lichens and bryophytes (de Graaf & Roberts 2009; Halpern et al. 2012; Rudolphi et al. 2014).”'
$d.Range($pStart + 205, $pStart + 1070).Font.Italic = $true

# --- paragraph 4: style=Normal ---
$para = $d.Paragraphs.Item($baseIndex + 5)
$para.Style = 'Normal'
$pStart = $para.Range.Start
$para.Range.Text = 'En nyligen publicerad vetenskaplig uppsats av Koelmeijer m.fl. (2022) inkluderar orkidén knärots skyddsbehov. I uppsatsen berörs problemet med uttorkning för växter, bl.a. för knärot, ett problem som blivit accentuerat på grund av den pågående klimatförändringen och torra somrar, t.ex. den exceptionellt torra sommaren 2018. I uppsatsen undersöks områden med tre olika avstånd från kalhyggeskant med avseende på skydd bl.a. för knärot. Det första området har avstånd upp till 20 m från hyggeskant (Strong edge effect), det andra 20 – 40 m från hyggeskant (Weak edge effect) och det tredje avser större avstånd från hyggeskant, där kanteffekten anses vara försumbar (Interior). Ett resultat var att man fann stor eller mycket stor uttorkningseffekt på känsliga och rödlistade skogsarter vid de kortare avstånden till hyggeskant, medan effekt av uttorkning inte konstaterades på större avstånd (Interior). För orkidén knärot fann man en rik förekomst (upp till 0,06 dm2/m2) på stort avstånd från hyggeskant (Interior), medan förekomsten var liten eller närmast försumbar i de områden som klassificerades som Weak edge effect respektive Strong edge effect. Arbetet påpekar att de allt oftare förekommande torra somrarna ger ytterligare skäl att utöka skyddsavståndet från hyggen till den fuktkrävande arten knärot (Koelmeijer m.fl., 2022).'

# --- paragraph 5: style=Normal ---
$para = $d.Paragraphs.Item($baseIndex + 6)
$para.Style = 'Normal'
$pStart = $para.Range.Start
$para.Range.Text = 'Även Skogsstyrelsens egen vägledning för hänsyn till knärot ligger i linje med ovanstående forskningsstudier. Av vägledningen framgår det att för med hög sannolikhet kunna bevara befintliga förekomster krävs relativt stora avsättningar av uppvuxen skog med slutet och relativt tätt kronskikt. Som riktlinje kan krävas ett avstånd på 50 meter in från brynet för att vidmakthålla ett fungerande mikroklimat. Detta innebär att fristående hänsynsytor för många arter (kärlväxter, lavar och mossor) kan behöva ha en area överstigande 0,8 hektar (cirkelyta med radien 50 meter = 0,78 hektar) för att bibehålla lokalklimatet. Även ganska små förändringar i form av förändrade ljus- och fuktighetsförhållanden, till exempel till följd av gallring, kan leda till att arten försvinner till följd av konkurrens med mera ljuskrävande och snabbväxande arter (Skogsstyrelsen, 2022).'

# --- paragraph 6: style=Heading2 ---
$para = $d.Paragraphs.Item($baseIndex + 7)
$para.Style = 'Heading2'
$pStart = $para.Range.Start
$para.Range.Text = 'Referenser - knärot'

# --- paragraph 7: style=Normal ---
$para = $d.Paragraphs.Item($baseIndex + 8)
$para.Style = 'Normal'
$pStart = $para.Range.Start
$para.Range.Text = 'de Graaf M & Roberts M.R., 2009. Short-term response of the herbaceous layer within leave patches after harvest. Forest Ecology and Management 257, 1014-1025'
$d.Range($pStart + 33, $pStart + 113).Font.Italic = $true

# --- paragraph 8: style=Normal ---
$para = $d.Paragraphs.Item($baseIndex + 9)
$para.Style = 'Normal'
$pStart = $para.Range.Start
$para.Range.Text = 'Halpern, C. B., Halaj, J., Evans, S. A., & Dovciak, M., 2012. Level and pattern of overstory retention interact to shape long-term responses of understories to timber harvest. Ecological Applications, 22, 2049-2064 '
$d.Range($pStart + 62, $pStart + 176).Font.Italic = $true

# --- paragraph 9: style=Normal ---
$para = $d.Paragraphs.Item($baseIndex + 10)
$para.Style = 'Normal'
$pStart = $para.Range.Start
$para.Range.Text = 'Koelmeijer, I. A., Ehrlén, J., Jönsson, M., De Frenne, P., Berg, P., Andersson, J., Weibull, H. & Hylander, N. 2022. Interactive effects of drought and edge exposure on old-growth forest understory species. Landscape Ecology, 37, sid 1839-1853'
$d.Range($pStart + 117, $pStart + 207).Font.Italic = $true

# --- paragraph 10: style=Normal ---
$para = $d.Paragraphs.Item($baseIndex + 11)
$para.Style = 'Normal'
$pStart = $para.Range.Start
$para.Range.Text = 'Rudolphi, J., Jönsson, M. T., & Gustafsson, L., 2014. Biological legacies buffer local species extinction after logging. Journal of Applied Ecology. 51, 53-62.'
$d.Range($pStart + 54, $pStart + 121).Font.Italic = $true

# --- paragraph 11: style=Normal ---
$para = $d.Paragraphs.Item($baseIndex + 12)
$para.Style = 'Normal'
$pStart = $para.Range.Start
$para.Range.Text = 'Skogsstyrelsen, 2022. Vägledning för hänsyn till knärot. https://www.skogsstyrelsen.se/lag-och-tillsyn/artskydd/vagledningar-och-kunskapsstod-artskydd/vagledning-for-hansyn-till-knarot/'
$d.Range($pStart + 22, $pStart + 57).Font.Italic = $true

# --- paragraph 12: style=Normal ---
$para = $d.Paragraphs.Item($baseIndex + 13)
$para.Style = 'Normal'
$pStart = $para.Range.Start
$para.Range.Text = 'SLU Artdatabanken, 2021. Artfaktablad. Naturvård – artfakta. SLU Artdatabanken, Uppsala '
$d.Range($pStart + 25, $pStart + 61).Font.Italic = $true

# --- update date in header ---
$sec = $d.Sections.First
for ($i = 1; $i -le $sec.Headers.Count; $i++) {
  $hdr = $sec.Headers.Item($i)
  if ($hdr.Exists) {
    $hdr.Range.Find.Execute("2023-09-13", $true, $false, $false, $false, $false, $true, 1, $false, "2023-09-15", 2)
  }
}
